$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.2233853333333333
$ws.Range("H2").Value = 0.670156
$ws.Range("I2").Value = 0.181392316683816
$ws.Range("J2").Value = 0.2059465991073217
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.009224
$ws.Range("N2").Value = 0.027672
$ws.Range("O2").Value = 0.01664709506025488
$ws.Range("P2").Value = 0.01697905168388481
$ws.Range("Q2").Value = 0.002060506314666667
$ws.Range("R2").Value = 0.018544556832
$ws.Range("S2").Value = 0.003019655139035341
$ws.Range("T2").Value = 0.003496777950363521
# Row 3
$ws.Range("G3").Value = 0.2233853333333333
$ws.Range("H3").Value = 0.670156
$ws.Range("I3").Value = 0.181392316683816
$ws.Range("J3").Value = 0.2059465991073217
$ws.Range("O3").Value = 0.9247000490894389
$ws.Range("P3").Value = 0.9431393206293146
$ws.Range("Q3").Value = 0.1144554220075556
$ws.Range("R3").Value = 1.030098798068
$ws.Range("S3").Value = 0.1677334841419717
$ws.Range("T3").Value = 0.1942363355679972
# Row 4
$ws.Range("G4").Value = 0.2233853333333333
$ws.Range("H4").Value = 0.670156
$ws.Range("I4").Value = 0.181392316683816
$ws.Range("J4").Value = 0.2059465991073217
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.5
$ws.Range("M4").Value = 0.032499
$ws.Range("N4").Value = 0.064998
$ws.Range("O4").Value = 0.05865285585030608
$ws.Range("P4").Value = 0.03988162768680055
$ws.Range("Q4").Value = 0.007259799948
$ws.Range("R4").Value = 0.043558799688
$ws.Range("S4").Value = 0.01063917740280893
$ws.Range("T4").Value = 0.008213485588960977
# Row 5
$ws.Range("I5").Value = 0.1777417664013226
$ws.Range("J5").Value = 0.2018018898423738
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.009224
$ws.Range("N5").Value = 0.027672
$ws.Range("O5").Value = 0.01664709506025488
$ws.Range("P5").Value = 0.01697905168388481
$ws.Range("Q5").Value = 0.002019038285333333
$ws.Range("R5").Value = 0.018171344568
$ws.Range("S5").Value = 0.002958884081460434
$ws.Range("T5").Value = 0.003426404717539294
# Row 6
$ws.Range("I6").Value = 0.1777417664013226
$ws.Range("J6").Value = 0.2018018898423738
$ws.Range("O6").Value = 0.9247000490894389
$ws.Range("P6").Value = 0.9431393206293146
$ws.Range("S6").Value = 0.1643578201165466
$ws.Range("T6").Value = 0.1903272972876482
# Row 7
$ws.Range("I7").Value = 0.1777417664013226
$ws.Range("J7").Value = 0.2018018898423738
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.5
$ws.Range("M7").Value = 0.032499
$ws.Range("N7").Value = 0.064998
$ws.Range("O7").Value = 0.05865285585030608
$ws.Range("P7").Value = 0.03988162768680055
$ws.Range("Q7").Value = 0.007113695276999999
$ws.Range("R7").Value = 0.04268217166199999
$ws.Range("S7").Value = 0.01042506220331555
$ws.Range("T7").Value = 0.008048187837186291
# Row 8
$ws.Range("G8").Value = 0.165029
$ws.Range("H8").Value = 0.4950869999999999
$ws.Range("I8").Value = 0.1340060790174831
$ws.Range("J8").Value = 0.1521458942578244
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.009224
$ws.Range("N8").Value = 0.027672
$ws.Range("O8").Value = 0.01664709506025488
$ws.Range("P8").Value = 0.01697905168388481
$ws.Range("Q8").Value = 0.001522227496
$ws.Range("R8").Value = 0.013700047464
$ws.Range("S8").Value = 0.002230811936056068
$ws.Range("T8").Value = 0.002583293002094474
# Row 9
$ws.Range("G9").Value = 0.165029
$ws.Range("H9").Value = 0.4950869999999999
$ws.Range("I9").Value = 0.1340060790174831
$ws.Range("J9").Value = 0.1521458942578244
$ws.Range("O9").Value = 0.9247000490894389
$ws.Range("P9").Value = 0.9431393206293146
$ws.Range("Q9").Value = 0.08455552366233332
$ws.Range("R9").Value = 0.760999712961
$ws.Range("S9").Value = 0.1239154278457499
$ws.Range("T9").Value = 0.1434947753468641
# Row 10
$ws.Range("G10").Value = 0.165029
$ws.Range("H10").Value = 0.4950869999999999
$ws.Range("I10").Value = 0.1340060790174831
$ws.Range("J10").Value = 0.1521458942578244
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.5
$ws.Range("M10").Value = 0.032499
$ws.Range("N10").Value = 0.064998
$ws.Range("O10").Value = 0.05865285585030608
$ws.Range("P10").Value = 0.03988162768680055
$ws.Range("Q10").Value = 0.005363277470999999
$ws.Range("R10").Value = 0.03217966482599999
$ws.Range("S10").Value = 0.007859839235677162
$ws.Range("T10").Value = 0.00606782590886588
# Row 11
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.5
$ws.Range("G11").Value = 0.4404835
$ws.Range("H11").Value = 0.8809669999999999
$ws.Range("I11").Value = 0.3576793576092536
$ws.Range("J11").Value = 0.2707312291105055
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.009224
$ws.Range("N11").Value = 0.027672
$ws.Range("O11").Value = 0.01664709506025488
$ws.Range("P11").Value = 0.01697905168388481
$ws.Range("Q11").Value = 0.004063019804
$ws.Range("R11").Value = 0.024378118824
$ws.Range("S11").Value = 0.005954322267212144
$ws.Range("T11").Value = 0.004596759531508932
# Row 12
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.5
$ws.Range("G12").Value = 0.4404835
$ws.Range("H12").Value = 0.8809669999999999
$ws.Range("I12").Value = 0.3576793576092536
$ws.Range("J12").Value = 0.2707312291105055
$ws.Range("O12").Value = 0.9247000490894389
$ws.Range("P12").Value = 0.9431393206293146
$ws.Range("Q12").Value = 0.2256895031001666
$ws.Range("R12").Value = 1.354137018601
$ws.Range("S12").Value = 0.3307461195395558
$ws.Range("T12").Value = 0.2553372674964215
# Row 13
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.5
$ws.Range("G13").Value = 0.4404835
$ws.Range("H13").Value = 0.8809669999999999
$ws.Range("I13").Value = 0.3576793576092536
$ws.Range("J13").Value = 0.2707312291105055
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.5
$ws.Range("M13").Value = 0.032499
$ws.Range("N13").Value = 0.064998
$ws.Range("O13").Value = 0.05865285585030608
$ws.Range("P13").Value = 0.03988162768680055
$ws.Range("Q13").Value = 0.0143152732665
$ws.Range("R13").Value = 0.057261093066
$ws.Range("S13").Value = 0.02097891580248563
$ws.Range("T13").Value = 0.01079720208257508
# Row 14
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.1837163333333333
$ws.Range("H14").Value = 0.551149
$ws.Range("I14").Value = 0.1491804802881247
$ws.Range("J14").Value = 0.1693743876819745
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.009224
$ws.Range("N14").Value = 0.027672
$ws.Range("O14").Value = 0.01664709506025488
$ws.Range("P14").Value = 0.01697905168388481
$ws.Range("Q14").Value = 0.001694599458666667
$ws.Range("R14").Value = 0.015251395128
$ws.Range("S14").Value = 0.002483421636490891
$ws.Range("T14").Value = 0.002875816482378587
# Row 15
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.1837163333333333
$ws.Range("H15").Value = 0.551149
$ws.Range("I15").Value = 0.1491804802881247
$ws.Range("J15").Value = 0.1693743876819745
$ws.Range("O15").Value = 0.9247000490894389
$ws.Range("P15").Value = 0.9431393206293146
$ws.Range("Q15").Value = 0.09413030903855556
$ws.Range("R15").Value = 0.847172781347
$ws.Range("S15").Value = 0.137947197445615
$ws.Range("T15").Value = 0.1597436449303835
# Row 16
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.1837163333333333
$ws.Range("H16").Value = 0.551149
$ws.Range("I16").Value = 0.1491804802881247
$ws.Range("J16").Value = 0.1693743876819745
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.5
$ws.Range("M16").Value = 0.032499
$ws.Range("N16").Value = 0.064998
$ws.Range("O16").Value = 0.05865285585030608
$ws.Range("P16").Value = 0.03988162768680055
$ws.Range("Q16").Value = 0.005970597117
$ws.Range("R16").Value = 0.035823582702
$ws.Range("S16").Value = 0.008749861206018807
$ws.Range("T16").Value = 0.006754926269212323
